# Added Week 15 simulations
$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 474
$wsOff.Range("C3").Value = 331
$wsOff.Range("D3").Value = 127
$wsOff.Range("E3").Value = 57
$wsOff.Range("F3").Value = 5
$wsOff.Range("G3").Value = 7

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 319
$wsDef.Range("C3").Value = 207
$wsDef.Range("D3").Value = 85
$wsDef.Range("E3").Value = 36
$wsDef.Range("F3").Value = 8
$wsDef.Range("G3").Value = 4
